$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J45").Value = 0.01855976243503714
$ws.Range("I46").Value = -0.1296176279974082
$ws.Range("H47").Value = -0.2870636170015632
$ws.Range("G48").Value = 0.2135958395245076
$ws.Range("F49").Value = -0.06676204101096155
$ws.Range("E50").Value = 0.1052128168340501
$ws.Range("D51").Value = -0.2006497229122814
$ws.Range("C52").Value = 0.4116802297750048
$ws.Range("B53").Value = -0.2766911554241067
